$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-6) ---

# Row 4: ending text changes, level/precondition stay the same (A4=3,B4=2,C4=1)
$ws.Range("D4").Value = "健身是一种信仰，你成为了一名健美教练。"

# Row 2: ending text changes, level/precondition stay the same (A2=1,B2=1,C2=0)
$ws.Range("D2").Value = "你活蹦乱跳，体能过人，在小学时体育老师发现了你的天赋，你立志成为一名运动员~"

# Row 3: level/precondition change, new ending text (A3=2,B3=1,C3=0)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "你对这五彩缤纷的世界有独特的美的感受，画画让你沉醉，成为一名艺术生快乐地画画是你的理想~"

# Row 5: level/precondition change, new ending text (A5=4,B5=2,C5=1)
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "你成为了一名国家级运动员，在赛场上挥洒汗水让你十分激动自豪。"

# Row 6: level/precondition change, new ending text (A6=5,B6=2,C6=1)
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "孩子们健康茁壮的成长让你欣慰，体育老师就是你热爱的职业。"

# --- Add new rows (7-11) for additional endings ---

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "体育也是一门科学，你考入大学进行体育理论相关的学习与深造。"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "艺术来自生活，你汲取生活中的美画成了多本漫画，作为一名漫画家你感到十分快乐。"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "你的画作在比赛上获奖，获得了众多赞美，被认可感与成就感让你更加热爱画画。"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "你的职业并不是画画相关，但画画让你的生活更加快乐，是你带给你诸多快乐的爱好。"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "你成为了一名美术老师，学生们可爱的画作让你感到自己工作的充实与幸福。"

# --- Update selection to match the post-edit active cell ---
$ws.Range("D15").Select()
